# Update NATMI ligand-receptor metrics for Ang-Plxnb2.xlsx with recomputed TPM values.
# Only the data rows (2-10) are affected; columns E,F,G,H are ligand-side metrics,
# M,N are receptor-side metrics, and I,J,O,P,Q,R,S,T are values derived from them
# (specificities / edge weights) that are recomputed downstream by the NATMI pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1825913333333333
$ws.Range("H2").Value = 0.547774
$ws.Range("I2").Value = 0.004164187522066756
$ws.Range("J2").Value = 0.004164187522066756
$ws.Range("M2").Value = 1.599392
$ws.Range("N2").Value = 4.798176
$ws.Range("O2").Value = 0.03952976301548796
$ws.Range("P2").Value = 0.03952976301548796
$ws.Range("Q2").Value = 0.2920351178026667
$ws.Range("R2").Value = 2.628316060224
$ws.Range("S2").Value = 0.0001646093458993509
$ws.Range("T2").Value = 0.0001646093458993509
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1825913333333333
$ws.Range("H3").Value = 0.547774
$ws.Range("I3").Value = 0.004164187522066756
$ws.Range("J3").Value = 0.004164187522066756
$ws.Range("O3").Value = 0.4638329693976876
$ws.Range("P3").Value = 0.4638329693976876
$ws.Range("Q3").Value = 3.426671589347556
$ws.Range("R3").Value = 30.840044304128
$ws.Range("S3").Value = 0.001931487463489022
$ws.Range("T3").Value = 0.001931487463489022
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1825913333333333
$ws.Range("H4").Value = 0.547774
$ws.Range("I4").Value = 0.004164187522066756
$ws.Range("J4").Value = 0.004164187522066756
$ws.Range("M4").Value = 20.09416733333333
$ws.Range("N4").Value = 60.28250199999999
$ws.Range("O4").Value = 0.4966372675868244
$ws.Range("P4").Value = 0.4966372675868245
$ws.Range("Q4").Value = 3.669020805616444
$ws.Range("R4").Value = 33.021187250548
$ws.Range("S4").Value = 0.002068090712678383
$ws.Range("T4").Value = 0.002068090712678383
$ws.Range("I5").Value = 0.943209744193024
$ws.Range("J5").Value = 0.9432097441930242
$ws.Range("M5").Value = 1.599392
$ws.Range("N5").Value = 4.798176
$ws.Range("O5").Value = 0.03952976301548796
$ws.Range("P5").Value = 0.03952976301548796
$ws.Range("Q5").Value = 66.14744588190933
$ws.Range("R5").Value = 595.327012937184
$ws.Range("S5").Value = 0.03728485766184926
$ws.Range("T5").Value = 0.03728485766184927
$ws.Range("I6").Value = 0.943209744193024
$ws.Range("J6").Value = 0.9432097441930242
$ws.Range("O6").Value = 0.4638329693976876
$ws.Range("P6").Value = 0.4638329693976876
$ws.Range("S6").Value = 0.4374917764138836
$ws.Range("T6").Value = 0.4374917764138837
$ws.Range("I7").Value = 0.943209744193024
$ws.Range("J7").Value = 0.9432097441930242
$ws.Range("M7").Value = 20.09416733333333
$ws.Range("N7").Value = 60.28250199999999
$ws.Range("O7").Value = 0.4966372675868244
$ws.Range("P7").Value = 0.4966372675868245
$ws.Range("Q7").Value = 831.0519536321908
$ws.Range("R7").Value = 7479.467582689717
$ws.Range("S7").Value = 0.4684331101172911
$ws.Range("T7").Value = 0.4684331101172912
$ws.Range("G8").Value = 2.307548333333334
$ws.Range("H8").Value = 6.922645
$ws.Range("I8").Value = 0.05262606828490914
$ws.Range("J8").Value = 0.05262606828490914
$ws.Range("M8").Value = 1.599392
$ws.Range("N8").Value = 4.798176
$ws.Range("O8").Value = 0.03952976301548796
$ws.Range("P8").Value = 0.03952976301548796
$ws.Range("Q8").Value = 3.690674343946667
$ws.Range("R8").Value = 33.21606909552
$ws.Range("S8").Value = 0.002080296007739345
$ws.Range("T8").Value = 0.002080296007739345
$ws.Range("G9").Value = 2.307548333333334
$ws.Range("H9").Value = 6.922645
$ws.Range("I9").Value = 0.05262606828490914
$ws.Range("J9").Value = 0.05262606828490914
$ws.Range("O9").Value = 0.4638329693976876
$ws.Range("P9").Value = 0.4638329693976876
$ws.Range("Q9").Value = 43.30550727971556
$ws.Range("R9").Value = 389.74956551744
$ws.Range("S9").Value = 0.02440970552031488
$ws.Range("T9").Value = 0.02440970552031488
$ws.Range("G10").Value = 2.307548333333334
$ws.Range("H10").Value = 6.922645
$ws.Range("I10").Value = 0.05262606828490914
$ws.Range("J10").Value = 0.05262606828490914
$ws.Range("M10").Value = 20.09416733333333
$ws.Range("N10").Value = 60.28250199999999
$ws.Range("O10").Value = 0.4966372675868244
$ws.Range("P10").Value = 0.4966372675868245
$ws.Range("Q10").Value = 46.36826233975444
$ws.Range("R10").Value = 417.31436105779
$ws.Range("S10").Value = 0.02613606675685492
$ws.Range("T10").Value = 0.02613606675685492
